# Rename the worksheet from "output-wrangling-2x2_monitor-in" to "Sheet1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# --- B2: Date_Seeded for FlyingM's first row changes from 7/25/2018 to 7/18/2018 ---
# Also needs to pick up the "corrected" (green highlight) formatting that matches C2.
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B2").Value = 43299         # 7/18/2018

# --- F2: PlotMix value changes from "Cool" to "Med-Warm" and loses its highlighted style ---
$ws.Range("D2").Copy()
$ws.Range("F2").PasteSpecial(-4122)   # xlPasteFormats (plain / no fill)
$ws.Range("F2").Value = "Med-Warm"

# --- Comments ---
# Remove the old threaded comment on F2 (no longer relevant now that F2 is fixed)
$ws.Range("F2").CommentThreaded.Delete()

# Add a new threaded comment on B2 explaining the Date_Seeded fix
$ws.Range("B2").AddCommentThreaded("FlyingM was seeded on 7/17 and 7/18. Every other observation for the 6/12/19 monitoring date has the Date_Seeded as 7/18.")

# --- Selection / active cell moved to F5 on save ---
$ws.Range("F5").Select()
